$wb = $excel.ActiveWorkbook

# --- "About" sheet: update unit-description text, drop the conversion-factor note row ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A11").Value = "The large primary energy output unit (used in totals graphs) is: quadrillion BTU"
$wsAbout.Range("A12").Value = "The small primary energy output unit (used in energy intensity per unit GDP graphs) is: thousand BTU"

# Row 15 ("1 Btu = 2.5219021687207...10-8 toe") is removed entirely.
$wsAbout.Range("A15").EntireRow.Delete()

# --- "BpTPEU-large" sheet: switch the large-unit conversion formula to 10^15 ---
$wsLarge = $wb.Worksheets.Item("BpTPEU-large")
$wsLarge.Range("B2").Formula = "=10^15"

# --- "BpTPEU-small" sheet: switch the small-unit conversion formula to 10^3 ---
$wsSmall = $wb.Worksheets.Item("BpTPEU-small")
$wsSmall.Range("B2").Formula = "=10^3"

# --- Restore default view state: reset stale selections, make "About" the active sheet ---
$wsLarge.Range("A1").Select()
$wsSmall.Range("A1").Select()
$wsAbout.Range("A1").Select()
$wsAbout.Activate()
